# Minor changes to SCMP, SQAP, and SPMP slides
#
# Slide 2 ("SPMP - Overview") — "Content Placeholder 2" shape:
#   - turn on PowerPoint's "shrink text on overflow" autofit (writes
#     <a:normAutofit/> into the shape's <a:bodyPr>) because the added
#     sentence below makes the paragraph overflow the placeholder
#   - prepend a sentence about the Team Leader's responsibility to the
#     third paragraph (which starts with "All ...")

$p = $ppt.ActivePresentation

$slide = $p.Slides.Item(2)
$contentShape = $slide.Shapes.Item(2)

# Enable "shrink text on overflow" (ppAutoSizeTextToFitShape) so the
# text frame picks up <a:normAutofit/> in its <a:bodyPr>.
$contentShape.TextFrame.AutoSize = 2

# Only the first run of paragraph 3 changes text ("All " -> longer
# sentence + "All "); the rest of the runs/paragraphs are untouched.
$contentShape.TextFrame.TextRange.Paragraphs(3).Runs(1).Text = `
    "The Team Leader is responsible for maintaining and updating this documents.  All "
